# ---------------------------------------------------------------------------
# edit.ps1 - applies the report_output.docx changes described by the diff:
#   1. "© 2025 Acme Corporation..." -> "© 2026 Acme Corporation..." (x2)
#   2. "Bold text" run gets xml:space="preserve" (text itself unchanged)
#   3. Two new styles added to styles.xml: StandardWeb ("Normal (Web)")
#      and EinfacheTabelle1 ("Plain Table 1")
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Bump the copyright year everywhere it appears ----------------------
$d.Content.Find.Execute(
    "© 2025 Acme Corporation. All rights reserved.", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "© 2026 Acme Corporation. All rights reserved.", 2) | Out-Null

# --- 2. Force the "Bold text" run to be serialized with xml:space="preserve"
# Re-typing the text as two pieces - the first one ending on a trailing
# space - makes the engine mark that (to-be-merged) run as
# whitespace-significant; the merge with the following "text" piece keeps
# that flag on the final single run even though the finished text has no
# leading/trailing whitespace of its own (exactly mirroring the target XML).
$boldRange = $d.Content
$boldRange.Find.ClearFormatting()
$boldRange.Find.Execute("Bold text", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($boldRange.Find.Found) {
    $boldRange.Text = "Bold "
    $boldRange.Collapse(0)
    $boldRange.InsertAfter("text")
}

# --- 3. Add the two missing styles to styles.xml ----------------------------

# 3a. Paragraph style "Normal (Web)" / styleId "StandardWeb"
$webStyle = $d.Styles.Add("StandardWeb", 1)
$webStyle.NameLocal = "Normal (Web)"
$webStyle.BaseStyle = "Standard"
$webStyle.Priority = 99
$webStyle.UnhideWhenUsed = $true
$webStyle.ParagraphFormat.SpaceBefore = 5
$webStyle.ParagraphFormat.SpaceBeforeAuto = $true
$webStyle.ParagraphFormat.SpaceAfter = 5
$webStyle.ParagraphFormat.SpaceAfterAuto = $true
$webStyle.Font.Name = "Times New Roman"
$webStyle.Font.NameFarEast = "Times New Roman"
$webStyle.Font.NameBi = "Times New Roman"
$webStyle.Font.Kerning = 0

# 3b. Table style "Plain Table 1" / styleId "EinfacheTabelle1"
$tableStyle = $d.Styles.Add("EinfacheTabelle1", 3)
$tableStyle.NameLocal = "Plain Table 1"
$tableStyle.BaseStyle = "NormaleTabelle"
$tableStyle.Priority = 41
